$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Layer0")
$ws1.Range("B2").Value = -0.7359652415508847
$ws1.Range("C2").Value = -0.4679439076167258
$ws1.Range("B3").Value = 0.6618991848898064
$ws1.Range("C3").Value = -1.651533521638236
$ws1.Range("B4").Value = -1.452238714845729
$ws1.Range("C4").Value = 0.2128442481488183

$ws2 = $wb.Worksheets.Item("Layer1")
$ws2.Range("B2").Value = -0.5515881728546594
$ws2.Range("C2").Value = -0.4692181449879241
$ws2.Range("B3").Value = -0.6965345639246838
$ws2.Range("C3").Value = 0.9101792829456896
$ws2.Range("B4").Value = -1.4515701270972
$ws2.Range("C4").Value = 0.8747949651556715
